# Add a new "as of" forecast date (2020-05-26) as column AG and row 45 to
# both the "cases" and "deaths" sheets, following the same staircase layout
# used by the rest of the table. Also backfill the newly-observed value for
# 2020-05-12 (row 31 / column B, the "Observed" column).

$wb = $excel.ActiveWorkbook

# Per-sheet data: new diagonal values for column AG (rows 32-45), and the
# newly observed value that belongs in B31.
$sheetData = @{
    "cases"  = @{
        B31 = 47719
        AG  = @{
            32 = 49043; 33 = 50184; 34 = 51143; 35 = 52320; 36 = 53397
            37 = 54496; 38 = 55485; 39 = 56380; 40 = 57514; 41 = 58283
            42 = 59026; 43 = 59821; 44 = 60630; 45 = 61332
        }
    }
    "deaths" = @{
        B31 = 3949
        AG  = @{
            32 = 4061; 33 = 4184; 34 = 4287; 35 = 4384; 36 = 4489
            37 = 4593; 38 = 4698; 39 = 4772; 40 = 4922; 41 = 5014
            42 = 5099; 43 = 5189; 44 = 5274; 45 = 5348
        }
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetData[$sheetName]

    # Materialize the new column AG (rows 2-44) as blank cells (matching the
    # existing blank cells under the other date columns) without pulling in
    # any new cell formatting.
    $ws.Range("AG2:AG44").Style = "Normal"

    # New column header AG1: reuses the existing "2020-05-12" label text.
    # Force text storage (quote-prefix) so it is NOT auto-converted to a
    # date serial number, matching the other header cells which are plain
    # text date labels.
    $ws.Range("AG1").Formula = "'2020-05-12"

    # New diagonal values for column AG.
    foreach ($r in $info.AG.Keys) {
        $ws.Cells.Item($r, 33).Value = $info.AG[$r]
    }

    # New row 45 for the "2020-05-26" as-of date.
    $ws.Range("B45:AF45").Style = "Normal"
    $ws.Range("A45").Formula = "'2020-05-26"

    # Newly observed value backfilled into the "Observed" column.
    $ws.Range("B31").Value = $info.B31
}
